$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Sheet1" to "Seats"
$ws.Name = "Seats"

# Move the active selection on the sheet from D3 to G22
$ws.Activate()
$ws.Range("G22").Select()
